# Store a reference to the worksheet, then stamp each cell (that has no
# value yet) with a label of its own address -- a quick diagnostic so the
# sheet shows, at a glance, which cell holds what.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: label every still-empty cell from row 3 down to row 14 with its
# own address (row 13 already holds a formula/value, so it is left as-is).
$ws.Range("B3").Value  = "B3"
$ws.Range("B4").Value  = "B4"
$ws.Range("B5").Value  = "B5"
$ws.Range("B6").Value  = "B6"
$ws.Range("B7").Value  = "B7"
$ws.Range("B8").Value  = "B8"
$ws.Range("B9").Value  = "B9"
$ws.Range("B10").Value = "B10"
$ws.Range("B11").Value = "B11"
$ws.Range("B12").Value = "B13"
$ws.Range("B14").Value = "B14"

# The one remaining empty cell in row 14 gets its own reference too.
$ws.Range("K14").Value = "K14"

# Leave the selection where the diagnostic run ended up.
$ws.Range("K15").Select() | Out-Null
